$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 "Non Recurring"
$ws.Range("D14").Value = -9034000
$ws.Range("E14").Value = -1994000

# Row 17 "Total Operating Expenses"
$ws.Range("D17").Value = -2219000
$ws.Range("E17").Value = 4408000

# Row 18 "Operating Income or Loss"
$ws.Range("D18").Value = 8879000
$ws.Range("E18").Value = 1690000

# Row 20 "Total Other Income/Expenses Net"
$ws.Range("D20").Value = 3046000
$ws.Range("E20").Value = 1000000

# Row 32 "Other Items"
$ws.Range("D32").Value = -3046000
$ws.Range("E32").Value = -1000000

# Row 91 "Capital Expenditures"
$ws.Range("D91").Value = -133000
$ws.Range("E91").Value = -152000
$ws.Range("F91").Value = -186000
$ws.Range("G91").Value = -262000
$ws.Range("H91").Value = -4235000
$ws.Range("I91").Value = -2439000
$ws.Range("J91").Value = -1837000
